# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# - Valor Mora total (E11): 764704 -> 1007056
# - Cant. Periodos (F13): 9 -> 10
# - Data table (rows 16..30) is re-sorted / extended (rows 16..34):
#     * LUISA CECILIA GUERRA POLO's 7 rows go from period desc (1701..1607)
#       to period asc (1607..1701)
#     * The remaining four workers (BRENDA, ANA, ROCIO, BELINDA) are
#       regrouped by period (2506, 2507, then a brand-new period 2508)
#       instead of being grouped by worker.
# - Signature rows shift from 35/36 down to 39/40 (handled automatically
#   by the row insert below).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Header summary cells
# ---------------------------------------------------------------------
$ws.Range("E11").Value = 1007056
$ws.Range("F13").Value = 10

# ---------------------------------------------------------------------
# 2. Make room for the new rows.
#    Before: data rows 16-30 (15 rows, last one styled as the table's
#    bottom border).
#    After:  data rows 16-34 (19 rows, last one keeps that bottom style).
#    Insert 4 blank rows right where the old row 30 lived so it is
#    pushed down to row 34, then copy the "normal" row format (row 29)
#    into the 4 freshly inserted rows (30-33).
# ---------------------------------------------------------------------
$ws.Rows.Item(30).Resize(4).Insert()

$ws.Range("B29:J29").Copy()
$ws.Range("B30:J33").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3. Fill in the data table, row by row.
# ---------------------------------------------------------------------
$tipoDoc = "CC"

$rows = @(
    @{ Doc = "33212762"; Nombre = "LUISA CECILIA GUERRA POLO"; Periodo = "1607"; Valor = 40000; Salario = 1200000 },
    @{ Doc = "33212762"; Nombre = "LUISA CECILIA GUERRA POLO"; Periodo = "1608"; Valor = 40000; Salario = 1200000 },
    @{ Doc = "33212762"; Nombre = "LUISA CECILIA GUERRA POLO"; Periodo = "1609"; Valor = 40000; Salario = 1200000 },
    @{ Doc = "33212762"; Nombre = "LUISA CECILIA GUERRA POLO"; Periodo = "1610"; Valor = 40000; Salario = 1200000 },
    @{ Doc = "33212762"; Nombre = "LUISA CECILIA GUERRA POLO"; Periodo = "1611"; Valor = 40000; Salario = 1200000 },
    @{ Doc = "33212762"; Nombre = "LUISA CECILIA GUERRA POLO"; Periodo = "1612"; Valor = 40000; Salario = 1200000 },
    @{ Doc = "33212762"; Nombre = "LUISA CECILIA GUERRA POLO"; Periodo = "1701"; Valor = 40000; Salario = 1200000 },

    @{ Doc = "33220302"; Nombre = "BRENDA ACENETH FLOREZ TURIZO";       Periodo = "2506"; Valor = 56940; Salario = 1423500 },
    @{ Doc = "33201157"; Nombre = "ANA ALFANIA ESCAÑO PIANETA";         Periodo = "2506"; Valor = 56940; Salario = 1423500 },
    @{ Doc = "51579311"; Nombre = "ROCIO DEL CARMEN BUSTAMANTE RODELO"; Periodo = "2506"; Valor = 71532; Salario = 1788305 },
    @{ Doc = "33216882"; Nombre = "BELINDA LEONOR HERRERA CASTRO";      Periodo = "2506"; Valor = 56940; Salario = 1423500 },

    @{ Doc = "33220302"; Nombre = "BRENDA ACENETH FLOREZ TURIZO";       Periodo = "2507"; Valor = 56940; Salario = 1423500 },
    @{ Doc = "33201157"; Nombre = "ANA ALFANIA ESCAÑO PIANETA";         Periodo = "2507"; Valor = 56940; Salario = 1423500 },
    @{ Doc = "51579311"; Nombre = "ROCIO DEL CARMEN BUSTAMANTE RODELO"; Periodo = "2507"; Valor = 71532; Salario = 1788305 },
    @{ Doc = "33216882"; Nombre = "BELINDA LEONOR HERRERA CASTRO";      Periodo = "2507"; Valor = 56940; Salario = 1423500 },

    @{ Doc = "33220302"; Nombre = "BRENDA ACENETH FLOREZ TURIZO";       Periodo = "2508"; Valor = 56940; Salario = 1423500 },
    @{ Doc = "33201157"; Nombre = "ANA ALFANIA ESCAÑO PIANETA";         Periodo = "2508"; Valor = 56940; Salario = 1423500 },
    @{ Doc = "51579311"; Nombre = "ROCIO DEL CARMEN BUSTAMANTE RODELO"; Periodo = "2508"; Valor = 71532; Salario = 1788305 },
    @{ Doc = "33216882"; Nombre = "BELINDA LEONOR HERRERA CASTRO";      Periodo = "2508"; Valor = 56940; Salario = 1423500 }
)

$r = 16
foreach ($row in $rows) {
    $ws.Cells.Item($r, 2).Value = $tipoDoc
    $ws.Cells.Item($r, 3).Value = $row.Doc
    $ws.Cells.Item($r, 4).Value = $row.Nombre
    $ws.Cells.Item($r, 5).Value = $row.Periodo
    $ws.Cells.Item($r, 6).Value = $row.Valor
    $ws.Cells.Item($r, 7).Value = $row.Salario
    $r++
}
